$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.370286
$ws.Range("H2").Value = 7.110858
$ws.Range("I2").Value = 0.3026841782318013
$ws.Range("J2").Value = 0.3026841782318014
$ws.Range("M2").Value = 46.29121633333333
$ws.Range("N2").Value = 138.873649
$ws.Range("O2").Value = 0.3133663986859022
$ws.Range("P2").Value = 0.3133663986859022
$ws.Range("Q2").Value = 109.7234219978713
$ws.Range("R2").Value = 987.5107979808421
$ws.Range("S2").Value = 0.09485105087170134
$ws.Range("T2").Value = 0.09485105087170136

$ws.Range("G3").Value = 2.370286
$ws.Range("H3").Value = 7.110858
$ws.Range("I3").Value = 0.3026841782318013
$ws.Range("J3").Value = 0.3026841782318014
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("N3").Value = 140.44867
$ws.Range("O3").Value = 0.3169204109998198
$ws.Range("P3").Value = 0.3169204109998198
$ws.Range("Q3").Value = 110.9678387398733
$ws.Range("R3").Value = 998.71054865886
$ws.Range("S3").Value = 0.09592679416836518
$ws.Range("T3").Value = 0.09592679416836519

$ws.Range("G4").Value = 2.370286
$ws.Range("H4").Value = 7.110858
$ws.Range("I4").Value = 0.3026841782318013
$ws.Range("J4").Value = 0.3026841782318014
$ws.Range("M4").Value = 38.53544233333333
$ws.Range("N4").Value = 115.606327
$ws.Range("O4").Value = 0.2608640200510233
$ws.Range("P4").Value = 0.2608640200510233
$ws.Range("Q4").Value = 91.34001946650733
$ws.Range("R4").Value = 822.060175198566
$ws.Range("S4").Value = 0.07895941153938814
$ws.Range("T4").Value = 0.07895941153938815

$ws.Range("G5").Value = 2.370286
$ws.Range("H5").Value = 7.110858
$ws.Range("I5").Value = 0.3026841782318013
$ws.Range("J5").Value = 0.3026841782318014
$ws.Range("M5").Value = 16.07945366666667
$ws.Range("N5").Value = 48.238361
$ws.Range("O5").Value = 0.1088491702632547
$ws.Range("P5").Value = 0.1088491702632547
$ws.Range("Q5").Value = 38.11290391374867
$ws.Range("R5").Value = 343.016135223738
$ws.Range("S5").Value = 0.03294692165234669
$ws.Range("T5").Value = 0.03294692165234669

$ws.Range("I6").Value = 0.2022126055089961
$ws.Range("J6").Value = 0.2022126055089961
$ws.Range("M6").Value = 46.29121633333333
$ws.Range("N6").Value = 138.873649
$ws.Range("O6").Value = 0.3133663986859022
$ws.Range("P6").Value = 0.3133663986859022
$ws.Range("Q6").Value = 73.3023416591041
$ws.Range("R6").Value = 659.721074931937
$ws.Range("S6").Value = 0.06336663595724715
$ws.Range("T6").Value = 0.06336663595724715

$ws.Range("I7").Value = 0.2022126055089961
$ws.Range("J7").Value = 0.2022126055089961
$ws.Range("M7").Value = 46.81622333333333
$ws.Range("N7").Value = 140.44867
$ws.Range("O7").Value = 0.3169204109998198
$ws.Range("P7").Value = 0.3169204109998198
$ws.Range("Q7").Value = 74.13369251863443
$ws.Range("R7").Value = 667.2032326677099
$ws.Range("S7").Value = 0.06408530204725546
$ws.Range("T7").Value = 0.06408530204725547

$ws.Range("I8").Value = 0.2022126055089961
$ws.Range("J8").Value = 0.2022126055089961
$ws.Range("M8").Value = 38.53544233333333
$ws.Range("N8").Value = 115.606327
$ws.Range("O8").Value = 0.2608640200510233
$ws.Range("P8").Value = 0.2608640200510233
$ws.Range("Q8").Value = 61.02103992175009
$ws.Range("R8").Value = 549.1893592957509
$ws.Range("S8").Value = 0.05274999317806843
$ws.Range("T8").Value = 0.05274999317806844

$ws.Range("I9").Value = 0.2022126055089961
$ws.Range("J9").Value = 0.2022126055089961
$ws.Range("M9").Value = 16.07945366666667
$ws.Range("N9").Value = 48.238361
$ws.Range("O9").Value = 0.1088491702632547
$ws.Range("P9").Value = 0.1088491702632547
$ws.Range("Q9").Value = 25.46188455879922
$ws.Range("R9").Value = 229.156961029193
$ws.Range("S9").Value = 0.02201067432642508
$ws.Range("T9").Value = 0.02201067432642508

$ws.Range("G10").Value = 2.286703333333333
$ws.Range("H10").Value = 6.860109999999999
$ws.Range("I10").Value = 0.2920107190904054
$ws.Range("J10").Value = 0.2920107190904054
$ws.Range("M10").Value = 46.29121633333333
$ws.Range("N10").Value = 138.873649
$ws.Range("O10").Value = 0.3133663986859022
$ws.Range("P10").Value = 0.3133663986859022
$ws.Range("Q10").Value = 105.8542786934878
$ws.Range("R10").Value = 952.6885082413899
$ws.Range("S10").Value = 0.09150634741904098
$ws.Range("T10").Value = 0.09150634741904098

$ws.Range("G11").Value = 2.286703333333333
$ws.Range("H11").Value = 6.860109999999999
$ws.Range("I11").Value = 0.2920107190904054
$ws.Range("J11").Value = 0.2920107190904054
$ws.Range("M11").Value = 46.81622333333333
$ws.Range("N11").Value = 140.44867
$ws.Range("O11").Value = 0.3169204109998198
$ws.Range("P11").Value = 0.3169204109998198
$ws.Range("Q11").Value = 107.0548139504111
$ws.Range("R11").Value = 963.4933255536998
$ws.Range("S11").Value = 0.09254415711048421
$ws.Range("T11").Value = 0.09254415711048421

$ws.Range("G12").Value = 2.286703333333333
$ws.Range("H12").Value = 6.860109999999999
$ws.Range("I12").Value = 0.2920107190904054
$ws.Range("J12").Value = 0.2920107190904054
$ws.Range("M12").Value = 38.53544233333333
$ws.Range("N12").Value = 115.606327
$ws.Range("O12").Value = 0.2608640200510233
$ws.Range("P12").Value = 0.2608640200510233
$ws.Range("Q12").Value = 88.11912443510776
$ws.Range("R12").Value = 793.0721199159698
$ws.Range("S12").Value = 0.07617509007991326
$ws.Range("T12").Value = 0.07617509007991326

$ws.Range("G13").Value = 2.286703333333333
$ws.Range("H13").Value = 6.860109999999999
$ws.Range("I13").Value = 0.2920107190904054
$ws.Range("J13").Value = 0.2920107190904054
$ws.Range("M13").Value = 16.07945366666667
$ws.Range("N13").Value = 48.238361
$ws.Range("O13").Value = 0.1088491702632547
$ws.Range("P13").Value = 0.1088491702632547
$ws.Range("Q13").Value = 36.76894029774555
$ws.Range("R13").Value = 330.9204626797099
$ws.Range("S13").Value = 0.03178512448096699
$ws.Range("T13").Value = 0.03178512448096699

$ws.Range("G14").Value = 1.590394666666667
$ws.Range("H14").Value = 4.771184
$ws.Range("I14").Value = 0.2030924971687972
$ws.Range("J14").Value = 0.2030924971687972
$ws.Range("M14").Value = 46.29121633333333
$ws.Range("N14").Value = 138.873649
$ws.Range("O14").Value = 0.3133663986859022
$ws.Range("P14").Value = 0.3133663986859022
$ws.Range("Q14").Value = 73.62130357004622
$ws.Range("R14").Value = 662.591732130416
$ws.Range("S14").Value = 0.06364236443791277
$ws.Range("T14").Value = 0.06364236443791277

$ws.Range("G15").Value = 1.590394666666667
$ws.Range("H15").Value = 4.771184
$ws.Range("I15").Value = 0.2030924971687972
$ws.Range("J15").Value = 0.2030924971687972
$ws.Range("M15").Value = 46.81622333333333
$ws.Range("N15").Value = 140.44867
$ws.Range("O15").Value = 0.3169204109998198
$ws.Range("P15").Value = 0.3169204109998198
$ws.Range("Q15").Value = 74.4562719028089
$ws.Range("R15").Value = 670.1064471252799
$ws.Range("S15").Value = 0.06436415767371494
$ws.Range("T15").Value = 0.06436415767371494

$ws.Range("G16").Value = 1.590394666666667
$ws.Range("H16").Value = 4.771184
$ws.Range("I16").Value = 0.2030924971687972
$ws.Range("J16").Value = 0.2030924971687972
$ws.Range("M16").Value = 38.53544233333333
$ws.Range("N16").Value = 115.606327
$ws.Range("O16").Value = 0.2608640200510233
$ws.Range("P16").Value = 0.2608640200510233
$ws.Range("Q16").Value = 61.28656196457421
$ws.Range("R16").Value = 551.5790576811679
$ws.Range("S16").Value = 0.05297952525365351
$ws.Range("T16").Value = 0.05297952525365351

$ws.Range("G17").Value = 1.590394666666667
$ws.Range("H17").Value = 4.771184
$ws.Range("I17").Value = 0.2030924971687972
$ws.Range("J17").Value = 0.2030924971687972
$ws.Range("M17").Value = 16.07945366666667
$ws.Range("N17").Value = 48.238361
$ws.Range("O17").Value = 0.1088491702632547
$ws.Range("P17").Value = 0.1088491702632547
$ws.Range("Q17").Value = 25.57267735438044
$ws.Range("R17").Value = 230.154096189424
$ws.Range("S17").Value = 0.02210644980351598
$ws.Range("T17").Value = 0.02210644980351598
